$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1 (section 5.2, first paragraph): the sentence was previously
# split across three runs ("...para que un " + "número" + " indeterminado
# ..."). Re-writing the whole sentence in one shot collapses it back into
# a single run while leaving the wording identical.
# ---------------------------------------------------------------------
$oldSentence1 = "El sistema estará diseñado para que un número indeterminado de cuentas de usuarios se registren, esta cantidad estará limitado a la capacidad de almacenamiento del hosting. En cuanto al flujo mensual de clientes registrados y sin registrar, se planea que en un principio sea de 500 usuarios, lo que equivale a un promedio de 16 personas al día. Por supuesto que, si esto no llega a ser suficiente, podría aumentarse el ancho de banda en el plan de hosting."

$rng = $d.Content
$rng.Find.ClearFormatting()
if ($rng.Find.Execute($oldSentence1, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)) {
    # Include the trailing paragraph mark so the whole paragraph is
    # rebuilt as a single run instead of leaving the original run split
    # untouched.
    $whole = $d.Range($rng.Start, $rng.End + 1)
    $whole.Text = $oldSentence1
}

# ---------------------------------------------------------------------
# Change 2 (section 5.2, "Al estar alojada..." paragraph): swap the
# sentence about the site running on any browser/OS for a note about the
# server itself running Linux.
# ---------------------------------------------------------------------
$oldSentence2 = "Al estar alojada en la web, la página funcionará correctamente en cualquier sistema operativo con el uso de cualquier navegador."
$newSentence2 = "El servidor estará alojado en Linux. "

$d.Content.Find.Execute($oldSentence2, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence2, 2)

# Drop the now-redundant blank paragraph left over in that part of
# section 5.2 -- it sits right before "REFERENCIAS", directly after the
# "Se contará con una base de datos..." paragraph, and carries no runs.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq [string][char]13) {
        $prev = $p.Previous()
        $next = $p.Next()
        if (($prev -ne $null) -and ($next -ne $null) -and `
            ($prev.Range.Text -like "*Se contará con una base de datos*") -and `
            ($next.Range.Text -like "*REFERENCIAS*")) {
            $p.Range.Delete()
            break
        }
    }
}
